$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text blocks for new/changed content ---
$objetivosText = @'
O conhecimento em assuntos que abordam os avanços na prevenção e controle da poluição. Discussões sobre a possibilidade de compensação das emissões, a inclusão de novas fontes, desde que protegidos os padrões de qualidade do ar, mostrar como algumas das melhores empresas mundiais estão aumentando a sua produtividade e os lucros com programas que também estão contribuindo para reduzir as emissões de poluentes são temas relevantes do futuro profissional. Demonstrar com exemplos práticos e as medidas tecnológicas adotadas em diversos setores da economia, como é possível aumentar a produtividade e a otimização dos lucros e, ainda, atingir a ecoeficiência.
'@

$resumidoText = @'
Apresentar os conceitos de ecoeficiência; tecnologias mais limpas; análise de ciclo de vida; desempenho ambiental; marketing ambiental; designer ambiental e contabilidade ambiental, aplicação e discussão de casos reais.
'@

$biblioText = @'
1)Biagio F. Giannetti, Cecília M. V. B. Almeida, Ecologia Industrial Conceitos, Ferramentas e Aplicações, 1ª Edição, Edgard Blucher, 2006, 128 p.
2)Joseph J. Romm, Empresas Eco-Eficientes, 1ª Edição, Signus Editora, 2004.
3)SALGADO, VIVIAN GULLO, Indicadores de ecoeficiência e o transporte de gás natural, 1ª Edição, Editora Interciencia, 2007, 117 p.
4)Canadian centre for pollution prevention (CCPP). Pollution Prevention Program Manual: P2 Planing and Beyond. Ontário (Canada): CCPP. 2001.
5)Canadian standards association. A guide to public involvement. Ontario (Canada): Etobicoke, 1996.
6)Sites: Association of chartered certified accountants (ACCA). Environmental, social and sustainability reporting on the world wide web: a guide to best practice. 2001 London: ACCA. Disponível em URL:http://www.accaglobal.com. Cowell, S. J. LCANET Theme Report: Positioning and Applications of LCA. Leiden (Holanda):Leiden University. 1997. Disponível em http://www.leidenuniv.nl/ Department for environment, food and rural affairs (DEFRA). Environmental reporting  general guidelines. London: DEFRA Publications. 2001. Disponível em http://www.defra.gov.uk. European commission. Integrated Pollution Prevention and Control (IPPC) Reference Document on Best Available Techniques. Sevilha (Espanha): European Commission, Joint Research Centre, Institute for Prospective Technological Studies Competitiveness and Sustainability Unit. 2001. Disponível em http://eippcb.jrc.es/reference. European environmental agency (EEA). Making Sustainability Accountable: Ecoefficiency, Resource Productivity and Innovation. In: Workshop on the fifth anniversary of the European environmental agency. 1998. Copenhague. Proceedings. Copenhague: EEA. 1998. Cleaner production. Cleaner production implementation. Copenhague: EEA: 2001. Disponível em URL:http://service.eea.int/envirowindows. Institute for global communications (IGC). Are business and industry taking sustainability seriously? San Francisco: IGC. 2001. Disponível em URL:http://www.igc.org/. International network for environmental management (INEM). The INEM sustainability reporting guide. Hamburg (Germany): INEM, 2001 Disponível em URL: http://www.inem.org/
7)Piotto, Z. C. Eco-eficiência na Indústria de Celulose e Papel - Estudo de Caso. 2003. 379 p. Tese (Doutorado) - Escola Politécnica, Universidade de São Paulo, São Paulo, 2003.
'@

# Row 10: B10/C10 objectives text changes (row height stays 60, unchanged)
$ws.Range("B10").Value = $objetivosText
$ws.Range("C10").Value = $objetivosText

# Row 13: remove A13 label; B13/C13 become the teacher name; row height -> default
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "2346890 - Eliane Corrêa Pedrozo"
$ws.Range("C13").Value = "2346890 - Eliane Corrêa Pedrozo"
$ws.Rows.Item(13).AutoFit()

# Row 14: A14 becomes "Programa resumido:"; B14/C14 new summary text; row height stays 60
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("B14").Value = $resumidoText
$ws.Range("C14").Value = $resumidoText
$ws.Rows.Item(14).RowHeight = 60

# Row 15: A15 becomes "Short syllabus:"; remove B15/C15; row height 120 -> 60
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Rows.Item(15).RowHeight = 60

# Row 16: A16 becomes "Programa:"; row height stays 120
$ws.Range("A16").Value = "Programa:"

# Row 17: A17 becomes "Syllabus:"; row height default -> 120
$ws.Range("A17").Value = "Syllabus:"
$ws.Rows.Item(17).RowHeight = 120

# Row 18: A18 becomes "Avaliação:"; remove B18/C18; row height 60 -> default
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18").Clear()
$ws.Range("C18").Clear()
$ws.Rows.Item(18).AutoFit()

# Row 19: A19 becomes "Método:"; B19/C19 unchanged; row height stays 60
$ws.Range("A19").Value = "Método:"

# Row 20: A20 becomes "Critério:"; B20/C20 unchanged; row height stays 60
$ws.Range("A20").Value = "Critério:"

# Row 21: A21 becomes "Norma de recuperação:"; B21/C21 unchanged; row height 120 -> 60
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Rows.Item(21).RowHeight = 60

# Row 22: A22 becomes "Bibliografia:"; B22/C22 new bibliography text; row height default -> 120
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B19").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("B22").Value = $biblioText
$ws.Range("C22").Value = $biblioText
$ws.Rows.Item(22).RowHeight = 120

# Row 23: remove B23/C23 (moved to row 24); add A23 "Requisitos:"; row height 30 -> default
$ws.Range("B23").Clear()
$ws.Range("C23").Clear()
$ws.Range("A22").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value = "Requisitos:"
$ws.Rows.Item(23).AutoFit()

# Row 24 (new row): B24/C24 get the prerequisite text; row height 30
$ws.Range("B19").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("B24").Value = "LOQ4073 -  Química Geral II  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOQ4073 -  Química Geral II  (Requisito fraco)`n"
$ws.Rows.Item(24).RowHeight = 30

Write-Output "Done"
